$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (hunk 0)
$ws.Range("H6").Value = 300
$ws.Range("I6").Value = 300
$ws.Range("K6").Value = 900
$ws.Range("M6").Value = -788

# Row 18 (hunk 1)
$ws.Range("H18").Value = 275000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# Row 40 (hunk 2)
$ws.Range("H40").Value = 4408.636
$ws.Range("I40").Value = 3648.75
$ws.Range("J40").Value = 4842.857
$ws.Range("K40").Value = 3648.75
$ws.Range("L40").Value = 4842.857
$ws.Range("M40").Value = -3473.75
$ws.Range("N40").Value = -5192.857

# Row 69 (hunk 3)
$ws.Range("H69").Value = 11752.167
$ws.Range("I69").Value = 7628.25
$ws.Range("K69").Value = 22884.75
$ws.Range("M69").Value = -22010.75

# Row 72 (hunk 4)
$ws.Range("H72").Value = 11752.167
$ws.Range("I72").Value = 7628.25
$ws.Range("K72").Value = 68654.25
$ws.Range("M72").Value = -64286.25

# Row 86 (hunk 5)
$ws.Range("H86").Value = 2635.3333
$ws.Range("I86").Value = 2571.8333
$ws.Range("K86").Value = 2571.8333
$ws.Range("M86").Value = -1448.8333

# Row 89 (hunk 6)
$ws.Range("H89").Value = 2635.3333
$ws.Range("I89").Value = 2571.8333
$ws.Range("K89").Value = 12859.1665
$ws.Range("M89").Value = -7243.166499999999

# Row 98 (hunk 7)
$ws.Range("H98").Value = 834.6
$ws.Range("I98").Value = 834.6
$ws.Range("K98").Value = 834.6
$ws.Range("M98").Value = 663.4

# Row 107 (hunk 8)
$ws.Range("H107").Value = 575.94446
$ws.Range("I107").Value = 575.94446
$ws.Range("K107").Value = 575.94446
$ws.Range("M107").Value = 1344.05554

# Row 122 (hunk 9)
$ws.Range("H122").Value = 834.6
$ws.Range("I122").Value = 834.6
$ws.Range("K122").Value = 2503.8
$ws.Range("M122").Value = -53.80000000000018

# Row 135 (hunk 10)
$ws.Range("H135").Value = 1077.3334
$ws.Range("J135").Value = 2210.5
$ws.Range("L135").Value = 19894.5
$ws.Range("N135").Value = -24964.5

# Row 137 (hunk 11)
$ws.Range("H137").Value = 3953.4546
$ws.Range("I137").Value = 1583.973
$ws.Range("K137").Value = 4751.919
$ws.Range("M137").Value = -2201.919

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk 12)
$ws.Range("H2").Value = 1838.95
$ws.Range("I2").Value = 1597.9667
$ws.Range("K2").Value = 1597.9667
$ws.Range("M2").Value = -1484.9667

# Row 74 (hunk 13)
$ws.Range("H74").Value = 1905.9375
$ws.Range("I74").Value = 1699.7333
$ws.Range("K74").Value = 1699.7333
$ws.Range("M74").Value = -825.7333000000001

# Row 77 (hunk 14)
$ws.Range("H77").Value = 1905.9375
$ws.Range("I77").Value = 1699.7333
$ws.Range("K77").Value = 8498.666500000001
$ws.Range("M77").Value = -4130.666500000001

# Row 116 (hunk 15)
$ws.Range("H116").Value = 1838.95
$ws.Range("I116").Value = 1597.9667
$ws.Range("K116").Value = 1597.9667
$ws.Range("M116").Value = 696.0333000000001

# Row 137 (hunk 16)
$ws.Range("H137").Value = 69256.664
$ws.Range("J137").Value = 69256.664
$ws.Range("L137").Value = 69256.664
$ws.Range("N137").Value = -79456.664

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk 17)
$ws.Range("H3").Value = 1838.95
$ws.Range("I3").Value = 1597.9667
$ws.Range("K3").Value = 1597.9667
$ws.Range("M3").Value = -1483.9667

# Row 58 (hunk 18)
$ws.Range("H58").Value = 65550.8
$ws.Range("J58").Value = 65550.8
$ws.Range("L58").Value = 65550.8
$ws.Range("N58").Value = -66138.8

# Row 86 (hunk 19)
$ws.Range("H86").Value = 2760.96
$ws.Range("I86").Value = 2202.125
$ws.Range("K86").Value = 2202.125
$ws.Range("M86").Value = -1079.125

# Row 89 (hunk 20)
$ws.Range("H89").Value = 2760.96
$ws.Range("I89").Value = 2202.125
$ws.Range("K89").Value = 11010.625
$ws.Range("M89").Value = -5394.625

# Row 99 (hunk 21)
$ws.Range("H99").Value = 2210.0588
$ws.Range("J99").Value = 5118.75
$ws.Range("L99").Value = 5118.75
$ws.Range("N99").Value = -8114.75

# Row 105 (hunk 22)
$ws.Range("H105").Value = 8813.362999999999
$ws.Range("I105").Value = 9518.5
$ws.Range("K105").Value = 9518.5
$ws.Range("M105").Value = -7771.5

# Row 129 (hunk 23)
$ws.Range("H129").Value = 66635
$ws.Range("J129").Value = 66635
$ws.Range("L129").Value = 66635
$ws.Range("N129").Value = -76635

# Row 134 (hunk 24)
$ws.Range("H134").Value = 1205.0741
$ws.Range("I134").Value = 1066.9131
$ws.Range("K134").Value = 3200.7393
$ws.Range("M134").Value = -665.7393000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (hunk 25)
$ws.Range("H16").Value = 1704.8334
$ws.Range("I16").Value = 1560.6154
$ws.Range("J16").Value = 2079.8
$ws.Range("K16").Value = 1560.6154
$ws.Range("L16").Value = 2079.8
$ws.Range("M16").Value = -1273.6154
$ws.Range("N16").Value = -2653.8

# Row 105 (hunk 26)
$ws.Range("H105").Value = 1847.1428
$ws.Range("J105").Value = 1854.2858
$ws.Range("L105").Value = 1854.2858
$ws.Range("N105").Value = -5348.2858

# Row 113 (hunk 27)
$ws.Range("H113").Value = 1704.8334
$ws.Range("I113").Value = 1560.6154
$ws.Range("J113").Value = 2079.8
$ws.Range("K113").Value = 1560.6154
$ws.Range("L113").Value = 2079.8
$ws.Range("M113").Value = 609.3846000000001
$ws.Range("N113").Value = -6419.8

# Row 119 (hunk 28)
$ws.Range("H119").Value = 44999.5
$ws.Range("J119").Value = 44999.5
$ws.Range("L119").Value = 44999.5
$ws.Range("N119").Value = -54675.5

# Row 134 (hunk 29)
$ws.Range("H134").Value = 5834.4443
$ws.Range("I134").Value = 5834.4443
$ws.Range("K134").Value = 17503.3329
$ws.Range("M134").Value = -14968.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (hunk 30)
$ws.Range("H2").Value = 41
$ws.Range("J2").Value = 41
$ws.Range("L2").Value = 246
$ws.Range("N2").Value = -472

# Row 7 (hunk 31)
$ws.Range("H7").Value = 89.916664
$ws.Range("I7").Value = 89
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 267
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -155
$ws.Range("N7").Value = -524

# Row 11 (hunk 32)
$ws.Range("H11").Value = 125541.75
$ws.Range("I11").Value = 562
$ws.Range("K11").Value = 1686
$ws.Range("M11").Value = -1546

# Row 22 (hunk 33)
$ws.Range("H22").Value = 490
$ws.Range("I22").Value = 490
$ws.Range("K22").Value = 1470
$ws.Range("M22").Value = -1301

# Row 27 (hunk 34)
$ws.Range("H27").Value = 490
$ws.Range("I27").Value = 490
$ws.Range("K27").Value = 1470
$ws.Range("M27").Value = -1368

# Row 39 (hunk 35)
$ws.Range("H39").Value = 4438.4443
$ws.Range("J39").Value = 5207
$ws.Range("L39").Value = 15621
$ws.Range("N39").Value = -16209

# Row 113 (hunk 36)
$ws.Range("H113").Value = 649.9167
$ws.Range("J113").Value = 644.3333
$ws.Range("L113").Value = 1932.9999
$ws.Range("N113").Value = -6272.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 14 (hunk 37)
$ws.Range("H14").Value = 15780
$ws.Range("I14").Value = 30700
$ws.Range("K14").Value = 30700
$ws.Range("M14").Value = -30532

# Row 123 (hunk 38)
$ws.Range("H123").Value = 67368.336
$ws.Range("J123").Value = 67368.336
$ws.Range("L123").Value = 67368.336
$ws.Range("N123").Value = -72268.336

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (hunk 39)
$ws.Range("H2").Value = 132.85715
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 132.85715
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 132.85715
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -356.85715

# Row 30 (hunk 40)
$ws.Range("H30").Value = 4129.5835
$ws.Range("I30").Value = 4129.5835
$ws.Range("K30").Value = 4129.5835
$ws.Range("M30").Value = -4021.5835

# Row 107 (hunk 41)
$ws.Range("H107").Value = 3999.5
$ws.Range("I107").Value = 3999.5
$ws.Range("K107").Value = 3999.5
$ws.Range("M107").Value = -2079.5

$ws = $wb.Worksheets.Item("WVR")
# Row 116 (hunk 42)
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 122 (hunk 43)
$ws.Range("H122").Value = 9604
$ws.Range("I122").Value = 9604
$ws.Range("K122").Value = 28812
$ws.Range("M122").Value = -26362

# Row 136 (hunk 44)
$ws.Range("H136").Value = 274.0909
$ws.Range("I136").Value = 274.0909
$ws.Range("K136").Value = 822.2727
$ws.Range("M136").Value = 1727.7273

# Row 138 (hunk 45)
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280
